$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change A4 from "RO.ACT.003HAB.SRL" to "AD.SEC.002.FON.01"
$ws.Range("A4").Value = "AD.SEC.002.FON.01"

# Add new rows in column D
$ws.Range("D8").Value = "RO.ACT.003HAB.SRA"
$ws.Range("D9").Value = "RO.ACT.003HAB.SRL"
$ws.Range("D10").Value = "RO.ACT.003HAB.SRM"
$ws.Range("D12").Value = "AD.SEC.002.FON.01"
$ws.Range("D14").Value = "RO.ACT.003MET.SRA"

# Change A3 from "RO.ACT.003HAB.SRA" to "RO.ACT.003MET"
$ws.Range("A3").Value = "RO.ACT.003MET"

# Update selection to A11
$ws.Range("A11").Select()
